# Update DCCD + 3 mm
# Recomputed "Neff" (column H) and "Nefferr" (column I) values for rows 2-25
# and 27-65 on Sheet1 (row 26 holds zeros and is intentionally left alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$block1 = New-Object 'object[,]' 24,2
$block1[0,0] = 151.5656107402997
$block1[0,1] = 0.147687521460849
$block1[1,0] = 152.005938202548
$block1[1,1] = 0.108366465094989
$block1[2,0] = 151.9683761157657
$block1[2,1] = 0.09319680967445086
$block1[3,0] = 152.0126328074481
$block1[3,1] = 0.1083533325169687
$block1[4,0] = 152.0115151722061
$block1[4,1] = 0.1161183228880781
$block1[5,0] = 154.0391088967869
$block1[5,1] = 0.2255241205133006
$block1[6,0] = 154.0524486596971
$block1[6,1] = 0.1809617531730773
$block1[7,0] = 155.2666129994101
$block1[7,1] = 0.2120535616338871
$block1[8,0] = 155.888396467415
$block1[8,1] = 0.2229883299108782
$block1[9,0] = 151.7750850039336
$block1[9,1] = 0.1649269456388452
$block1[10,0] = 150.8166232066449
$block1[10,1] = 0.1755595128831169
$block1[11,0] = 151.3417261093441
$block1[11,1] = 0.1252485475220408
$block1[12,0] = 151.3656736527939
$block1[12,1] = 0.1205553893788476
$block1[13,0] = 152.1741395845495
$block1[13,1] = 0.1309390400281521
$block1[14,0] = 155.7401109774144
$block1[14,1] = 0.1888335954364468
$block1[15,0] = 152.0699525618234
$block1[15,1] = 0.1035737638768537
$block1[16,0] = 152.0937504175068
$block1[16,1] = 0.1230403917271251
$block1[17,0] = 152.2840680321257
$block1[17,1] = 0.09519074387375316
$block1[18,0] = 152.0328472766903
$block1[18,1] = 0.1232287693425566
$block1[19,0] = 157.7916612115569
$block1[19,1] = 0.1094551823307176
$block1[20,0] = 151.9414710742979
$block1[20,1] = 0.104541251533543
$block1[21,0] = 152.0897904988037
$block1[21,1] = 0.1025951813059523
$block1[22,0] = 152.4628663311817
$block1[22,1] = 0.1137220878179702
$block1[23,0] = 152.154778507928
$block1[23,1] = 0.1023477671918774
$ws.Range("H2:I25").Value2 = $block1

$block2 = New-Object 'object[,]' 39,2
$block2[0,0] = 151.9819786895464
$block2[0,1] = 0.1212146117422016
$block2[1,0] = 152.0467195731177
$block2[1,1] = 0.1416509681458803
$block2[2,0] = 151.8715862216874
$block2[2,1] = 0.105514303983595
$block2[3,0] = 151.9778985141396
$block2[3,1] = 0.1251530014848004
$block2[4,0] = 151.6539673010244
$block2[4,1] = 0.1214680593578193
$block2[5,0] = 151.1502515658681
$block2[5,1] = 0.1322435614617539
$block2[6,0] = 151.7116636126021
$block2[6,1] = 0.1107774723696583
$block2[7,0] = 151.4815844843498
$block2[7,1] = 0.1252470967066953
$block2[8,0] = 152.0470186025359
$block2[8,1] = 0.1177096401923902
$block2[9,0] = 152.1316784320978
$block2[9,1] = 0.1227186168363573
$block2[10,0] = 152.2975699134284
$block2[10,1] = 0.1070752069472101
$block2[11,0] = 152.0934009190962
$block2[11,1] = 0.1177086764243578
$block2[12,0] = 152.4949482758668
$block2[12,1] = 0.1405104642107665
$block2[13,0] = 152.2190138412069
$block2[13,1] = 0.1266675211688783
$block2[14,0] = 152.5203401905414
$block2[14,1] = 0.1357392528821232
$block2[15,0] = 152.0640662767559
$block2[15,1] = 0.1387604532006805
$block2[16,0] = 152.08846041115
$block2[16,1] = 0.1070804885658267
$block2[17,0] = 152.0563785797918
$block2[17,1] = 0.1203537104700684
$block2[18,0] = 152.2735879201218
$block2[18,1] = 0.1116148848110251
$block2[19,0] = 152.242116198506
$block2[19,1] = 0.1266673590295916
$block2[20,0] = 151.932216149823
$block2[20,1] = 0.1267344387142841
$block2[21,0] = 152.3701687694523
$block2[21,1] = 0.1305112869875901
$block2[22,0] = 152.1737753034404
$block2[22,1] = 0.1234163987427716
$block2[23,0] = 152.2009649183906
$block2[23,1] = 0.1266684343392947
$block2[24,0] = 152.1200200937825
$block2[24,1] = 0.1091827579426642
$block2[25,0] = 152.1674935914248
$block2[25,1] = 0.1064786397919329
$block2[26,0] = 151.9270346561076
$block2[26,1] = 0.1143448201333478
$block2[27,0] = 151.9667388777131
$block2[27,1] = 0.1233654465559598
$block2[28,0] = 152.0716831063736
$block2[28,1] = 0.1109824304626085
$block2[29,0] = 152.0572976558605
$block2[29,1] = 0.1091846526434656
$block2[30,0] = 152.0165788028601
$block2[30,1] = 0.1166457401855792
$block2[31,0] = 152.0115031615022
$block2[31,1] = 0.1233645712220186
$block2[32,0] = 151.8284493428507
$block2[32,1] = 0.1177143689433772
$block2[33,0] = 151.9510377803976
$block2[33,1] = 0.117711815773257
$block2[34,0] = 151.5376930428534
$block2[34,1] = 0.1216482693420094
$block2[35,0] = 152.0342926674716
$block2[35,1] = 0.1114586980150636
$block2[36,0] = 152.0355160329138
$block2[36,1] = 0.1050548728597248
$block2[37,0] = 152.1246091290602
$block2[37,1] = 0.121087109404861
$block2[38,0] = 152.0353634556942
$block2[38,1] = 0.1205504444102803
$ws.Range("H27:I65").Value2 = $block2
